$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Resize the workbook window (matches the larger window the file was
# last saved from).
try {
    $win = $wb.Windows.Item(1)
    $win.Width = 19200
    $win.Height = 7140
} catch {}

# Clear the sample/demo data row (row 2, columns A:F) that shipped with the template.
$ws.Range("A2:F2").ClearContents()

# Move the active selection, matching how the file was left after editing.
$ws.Range("B5").Select()
